# Update odds values in row 7 of Sheet1 as per the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "G7"  = 5.4
    "H7"  = 3.75
    "K7"  = 2.18
    "L7"  = 2.1
    "T7"  = 2.72
    "U7"  = 1.95
    "V7"  = 1.75
    "W7"  = 13.5
    "X7"  = 32
    "Y7"  = 18
    "Z7"  = 110
    "AA7" = 65
    "AB7" = 65
    "AD7" = 7.5
    "AE7" = 18.5
    "AF7" = 100
    "AG7" = 6.3
    "AI7" = 8.25
    "AJ7" = 10.75
    "AM7" = 900
    "AN7" = 7
    "AO7" = 35
    "AP7" = 40
    "AQ7" = 250
    "AR7" = 300
    "AT7" = 2.72
    "AU7" = 8.25
    "AV7" = 90
    "AW7" = 3.3
    "AX7" = 7.4
    "AY7" = 18
    "AZ7" = 24
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
